# Apply the "DatesBackTo" / "Localimage" popup corrections to the Geo sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geo")

# New "Attribution" header column.
$ws.Range("I1").Value = "Attribution"

# DatesBackTo (column E) becomes text like "1350 CE" instead of a bare number.
$ws.Range("E2").Value = "1350 CE"
$ws.Range("E3").Value = "1526 CE"
$ws.Range("E4").Value = "1600 CE"
$ws.Range("E5").Value = "1600 CE"
$ws.Range("E6").Value = "1630 CE"
$ws.Range("E7").Value = "1200 CE"
$ws.Range("E8").Value = "1998 CE"
$ws.Range("E9").Value = "1408 CE"
$ws.Range("E10").Value = "1503 CE"
$ws.Range("E11").Value = "1529 CE"
$ws.Range("E12").Value = "1619 CE"

# Replace the weird Vefsn harpa image with a cleaner local asset.
$ws.Range("H4").Value = "images/vefsen.png"

# Update the active selection to reflect where the author left off.
$ws.Range("H5").Select() | Out-Null
